# Creando pagina de exel de cliente
# Adds a new "Client" worksheet (as the last sheet) with a header row:
# Doc Id | Name | Address | Contac Number | Stratum | Description

$wb = $excel.ActiveWorkbook

# Insert the new sheet after the current last sheet so it lands at the end.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "Client"

$ws.Range("A1").Value = "Doc Id"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Address"
$ws.Range("D1").Value = "Contac Number"
$ws.Range("E1").Value = "Stratum"
$ws.Range("F1").Value = "Description"

# Keep the originally active sheet ("Animal") selected/active, as before.
$wb.Worksheets.Item("Animal").Activate()
